$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

# 1. Title heading "Asignación de Rol" -> "Asociación de tarjeta".
# The phrase "Asignación de Rol" also occurs later as part of a different
# heading ("Caso de uso: Asignación de Rol") that must stay untouched, so we
# scope the Find/Replace to the specific paragraph whose whole text is the
# title, instead of doing a document-wide replace.
$titlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Asignación de Rol`r") {
        $titlePara = $p
        break
    }
}
if ($titlePara -eq $null) {
    throw "Title paragraph 'Asignación de Rol' not found"
}
$ok = $titlePara.Range.Find.Execute("Asignación de Rol", $true, $false, $false, $false, $false, $true, 1, $false, "Asociación de tarjeta", 2)
if (-not $ok) {
    throw "Failed to replace title text"
}

# 2. Short description paragraph
Replace-Text "usuario asignar su rol de worker." "usuario asociar su tarjeta a la aplicación."

# 3. "está logeado." text is unchanged, but the source merges the runs and
#    drops the spell-check proofErr markers around "logeado". Re-applying a
#    same-text replace triggers the same run-merge/proofErr-removal here.
Replace-Text "está logeado." "está logeado."

# 4. Precondition: Use case name reference
Replace-Text "del Caso de Uso: Actualizar Datos de Cuenta" "del Caso de Uso: Asociar tarjeta"

# 5. "hace click en el botón Ofrecer Servicios." text is unchanged, but the
#    source merges the runs and drops the proofErr markers around "click".
Replace-Text "hace click en el botón Ofrecer Servicios." "hace click en el botón Ofrecer Servicios."

# 6. Flow table row 2 - system response
Replace-Text "despliega un listado que contiene todos los servicios que se pueden ofrecer en la plataforma. A su vez, mostrará botones de Guardar y Cancelar." "despliega la vista de asociación de tarjeta, donde le presenta un formulario a rellenar."

# 7. Flow table row 3 - actor action
Replace-Text "selecciona los servicios que ofrecerá y hará click en guardar. En el caso que seleccione cancelar se activara el flujo alternativo 2.2.1" "completa los campos obligatorios del formulario y acepta para finalizar la asociación."

# 8. Flow table row 4 - system updates database
Replace-Text "actualiza la base de datos cambiando el rol del usuario de cliente a worker." "actualiza la tabla correspondiente en la database."

# 9. Postcondition
Replace-Text "se convierte en worker." "puede pagar un servicio sin necesidad de volver a completar el formulario de información de su tarjeta."
